$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "DAPI_ratio" in J1
$ws.Range("J1").Value = "DAPI_ratio"

# Add the per-row DAPI ratio formulas (mirrors the pattern used for the
# other "ratio" helper columns already present in the sheet)
$ws.Range("J2").Formula = "=E2/B2"
$ws.Range("J3").Formula = "=B3/E3"
$ws.Range("J4").Formula = "=B4/E4"

# Update the active cell/selection left behind after the edit
$ws.Range("H10").Select()
